# Deploying to gh-pages — add the 2022 column (U) to the sheet, mirroring
# the existing 2021 column (T) for formatting, then fill in the new figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column T's formatting (rows 3-40, the bordered data block) into the
# new column U so the new year column matches the existing style (number
# formats, borders, fills, etc.)
$ws.Range("T3:T40").Copy() | Out-Null
$ws.Range("U3:U40").PasteSpecial(-4122) | Out-Null

# Row 4 header: new year value 2022
$ws.Range("U4").Value = 2022

# Section "Boys" (A5 header) totals and causes, 2022 column
$ws.Range("U6").Value = 1456
$ws.Range("U8").Value = 45
$ws.Range("U9").Value = 35
$ws.Range("U10").Value = "-"
$ws.Range("U11").Value = 217
$ws.Range("U12").Value = 22
$ws.Range("U13").Value = 8
$ws.Range("U14").Value = "-"
$ws.Range("U15").Value = "-"
$ws.Range("U16").Value = 57
$ws.Range("U17").Value = "-"
$ws.Range("U18").Value = 5
$ws.Range("U19").Value = "-"
$ws.Range("U20").Value = 46
$ws.Range("U21").Value = 1021
$ws.Range("U22").Value = "-"

# Section "Girls" (A23 header) totals and causes, 2022 column
$ws.Range("U24").Value = 1019
$ws.Range("U26").Value = 15
$ws.Range("U27").Value = 30
$ws.Range("U28").Value = 1
$ws.Range("U29").Value = 179
$ws.Range("U30").Value = 16
$ws.Range("U31").Value = 8
$ws.Range("U32").Value = "-"
$ws.Range("U33").Value = "-"
$ws.Range("U34").Value = 46
$ws.Range("U35").Value = "-"
$ws.Range("U36").Value = "-"
$ws.Range("U37").Value = "-"
$ws.Range("U38").Value = 25
$ws.Range("U39").Value = 699
$ws.Range("U40").Value = "-"

# Move the selection/view to where the author left it after adding the column
$ws.Range("V6").Select() | Out-Null
